# Commit: "Add files via upload"
#
# The "Classes" sheet (tab name "Classes", internal sheet2.xml) gets two new
# trailing columns, X and Y, added right after the existing last column W
# (rows 1-43):
#   - X1 = "CategoriaRvt", Y1 = "ClasseIfc"   (new header labels)
#   - X2:X43 and Y2:Y43 = "null"              (placeholder values, same
#     shared-string "null" used elsewhere on this sheet, e.g. column G-K)
#
# Formatting for the new cells mirrors the existing column V ("Fonte"-like
# column immediately to the left of the computed Key column W): the header
# cell style from V1 and the body-cell style from V2:V43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")

# --- Header row (row 1): copy style from V1, then set the two new labels ---
$null = $ws.Range("V1").Copy()
$null = $ws.Range("X1:Y1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("X1").Value = "CategoriaRvt"
$ws.Range("Y1").Value = "ClasseIfc"

# --- Data rows (2-43): copy style from V2:V43, then fill with "null" ---
$null = $ws.Range("V2:V43").Copy()
$null = $ws.Range("X2:X43").PasteSpecial(-4122)  # xlPasteFormats
$null = $ws.Range("Y2:Y43").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 24).Value = "null"
    $ws.Cells.Item($r, 25).Value = "null"
}

# Match the author's final selection, now resting on the newly added columns.
$null = $ws.Range("X2:Y43").Select()
